$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")

# Row 32's "No." column had been typed as =A30+1, breaking the otherwise
# consistent +1 sequence shared by the surrounding rows (A29:A39). Fix it
# to follow the same pattern as its neighbors; this also renumbers the
# cascading rows below it (A33:A39) on recalculation.
$ws.Range("A32").Formula = "=A31+1"

# Leave the same cell selected that was active when the sheet was last saved.
$ws.Activate()
$null = $ws.Range("A35").Select()
